# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp in A1
# - A handful of countries changed rank (total cases moved them past a
#   neighbouring row), so the country label + stats for those specific
#   rows need updating in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Junio de 2020 a las 06:37"

# --- rank swap: Mexico / Pakistan (rows 16-17) ---
$ws.Cells.Item(16, 1).Value = "Pakistan"
$ws.Cells.Item(16, 2).Value = 181088
$ws.Cells.Item(16, 3).Value = 4471
$ws.Cells.Item(16, 4).Value = 71458
$ws.Cells.Item(16, 5).Value = 106040
$ws.Cells.Item(16, 7).Value = 89
$ws.Cells.Item(16, 8).Value = 3590

$ws.Cells.Item(17, 1).Value = "Mexico"
$ws.Cells.Item(17, 2).Value = 180545
$ws.Cells.Item(17, 3).Value = 5343
$ws.Cells.Item(17, 4).Value = 135279
$ws.Cells.Item(17, 5).Value = 23441
$ws.Cells.Item(17, 7).Value = 1044
$ws.Cells.Item(17, 8).Value = 21825

# --- rank swap: Austria / Kazajistan (rows 55-56) ---
$ws.Cells.Item(55, 1).Value = "Kazajistan"
$ws.Cells.Item(55, 2).Value = 17732
$ws.Cells.Item(55, 3).Value = 507
$ws.Cells.Item(55, 4).Value = 10897
$ws.Cells.Item(55, 5).Value = 6715
$ws.Cells.Item(55, 8).Value = 120

$ws.Cells.Item(56, 1).Value = "Austria"
$ws.Cells.Item(56, 2).Value = 17341
$ws.Cells.Item(56, 4).Value = 16197
$ws.Cells.Item(56, 5).Value = 454
$ws.Cells.Item(56, 8).Value = 690

# --- rank rotation: Bulgaria / Venezuela / Bosnia y Herzegovina / Grecia / Kirguistan (rows 90-94) ---
$ws.Cells.Item(90, 1).Value = "Venezuela"
$ws.Cells.Item(90, 2).Value = 3917
$ws.Cells.Item(90, 4).Value = 835
$ws.Cells.Item(90, 5).Value = 3049
$ws.Cells.Item(90, 8).Value = 33

$ws.Cells.Item(91, 1).Value = "Bulgaria"
$ws.Cells.Item(91, 2).Value = 3905
$ws.Cells.Item(91, 4).Value = 2074
$ws.Cells.Item(91, 5).Value = 1632
$ws.Cells.Item(91, 8).Value = 199

$ws.Cells.Item(92, 1).Value = "Kirguistan"
$ws.Cells.Item(92, 2).Value = 3356
$ws.Cells.Item(92, 3).Value = 205
$ws.Cells.Item(92, 4).Value = 2021
$ws.Cells.Item(92, 5).Value = 1295
$ws.Cells.Item(92, 7).Value = 3
$ws.Cells.Item(92, 8).Value = 40

$ws.Cells.Item(93, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(93, 2).Value = 3273
$ws.Cells.Item(93, 4).Value = 2241
$ws.Cells.Item(93, 5).Value = 863
$ws.Cells.Item(93, 8).Value = 169

$ws.Cells.Item(94, 1).Value = "Grecia"
$ws.Cells.Item(94, 2).Value = 3266
$ws.Cells.Item(94, 4).Value = 1374
$ws.Cells.Item(94, 5).Value = 1702
$ws.Cells.Item(94, 8).Value = 190

# --- data refresh only (no rank change): Mongolia (row 164) ---
$ws.Cells.Item(164, 2).Value = 213
$ws.Cells.Item(164, 3).Value = 7
$ws.Cells.Item(164, 4).Value = 153
$ws.Cells.Item(164, 5).Value = 60

# --- rank swap: Islas Virgenes Britanicas / Papua Nueva Guinea (rows 213-214) ---
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
